$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the Resources table: split the old combined "Notes" column
# --- into separate "Analog to" (C), "Unit" (D) and "Notes" (E) columns, and add
# --- untradeable resources (farm/factory) with their own units.

# Row 1
$ws.Range("A1").Value = "Resources"
$ws.Range("B1").Value = "Weight"
$ws.Range("C1").Value = "Analog to"
$ws.Range("D1").Value = "Unit"
$ws.Range("E1").Value = "Notes"

# Row 2
$ws.Range("A2").Value = "R1"
$ws.Range("B2").Value = 0.2
$ws.Range("C2").Value = "population"
$ws.Range("D2").Value = "million people"
$ws.Range("E2").Value = " the amount of people in a country is only a small indicator of the country's prosperity. Some small countries are very wealthy and some are very poor. Overall, does indicate some sense of wealth for a country."

# Row 3
$ws.Range("A3").Value = "R2"
$ws.Range("B3").Value = 0.75
$ws.Range("C3").Value = "metallic elements"
$ws.Range("D3").Value = "million tons"
$ws.Range("E3").Value = "essential for metallic alloy creation and electronic creation "

# Row 4
$ws.Range("A4").Value = "R3"
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = "timber"
$ws.Range("D4").Value = "million tons"
$ws.Range("E4").Value = "used in all forms of construction, but not a particularly rare resource"

# Row 5
$ws.Range("A5").Value = "R4"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "available land"
$ws.Range("D5").Value = "million acres"
$ws.Range("E5").Value = "valued at twice the weight as water because land limits how much housing/farm/factories can be created that bring large amounts of prosperity. "

# Row 6
$ws.Range("A6").Value = "R5"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = "renewable energy"
$ws.Range("D6").Value = "million kW"
$ws.Range("E6").Value = "renewable energy valued at 1 - in direct correlation to renewable energy waste's weight being -1"

# Row 7
$ws.Range("A7").Value = "R6"
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = "fossil fuel energy"
$ws.Range("D7").Value = "million kW"
$ws.Range("E7").Value = "fossil fuels create more energy than green sources, but their waste is higher to indicate penalty for using nonrenewables."

# Row 8
$ws.Range("A8").Value = "R7"
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = "water"
$ws.Range("D8").Value = "billion gallons"
$ws.Range("E8").Value = "0.5 chosen as the baseline for which all other raw resources are weighted. Essential for life and is involved in other types of resource creation, but is not rare."

# Row 9
$ws.Range("A9").Value = "R8"
$ws.Range("B9").Value = 0.5
$ws.Range("C9").Value = "animals"
$ws.Range("D9").Value = "million animals"
$ws.Range("E9").Value = "used for farms and food. Not particularly rare and has only a few use cases"

# Row 10
$ws.Range("A10").Value = "R9"
$ws.Range("B10").Value = 0.5
$ws.Range("C10").Value = "plants"
$ws.Range("D10").Value = "million tons"
$ws.Range("E10").Value = "used for farms and food, also produces fresh oxygen. Not particularly rare as well."

# Row 11
$ws.Range("A11").Value = "R18"
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "metallic alloys"
$ws.Range("D11").Value = "million tons"
$ws.Range("E11").Value = "weighted at 2 to account for -1 alloy waste weight. Alloy + alloy waste = 2 in weight, compared to 1.5 in lost input resources"

# Row 12
$ws.Range("A12").Value = "R19"
$ws.Range("B12").Value = 15
$ws.Range("C12").Value = "housing"
$ws.Range("D12").Value = "million homes"
$ws.Range("E12").Value = "weighted at 15 to account for -2 housing waste weight. Input resources lost have combined weight of 12.25."

# Row 13
$ws.Range("A13").Value = "R20"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "electronics"
$ws.Range("D13").Value = "million gadgets"
$ws.Range("E13").Value = "weighted at 5 since 2 electronics and 1 waste is created - these total to 9 in weight compared to 8.25 of lost input resources"

# Row 14
$ws.Range("A14").Value = "R21"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "farm"
$ws.Range("D14").Value = "million acres"

# Row 15
$ws.Range("A15").Value = "R22"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "factory"
$ws.Range("D15").Value = "thousand factories"

# Row 16
$ws.Range("A16").Value = "R1'"
$ws.Range("B16").Value = -1
$ws.Range("C16").Value = "population waste"

# Row 17
$ws.Range("A17").Value = "R5'"
$ws.Range("B17").Value = -1
$ws.Range("C17").Value = "renewable energy waste"
$ws.Range("E17").Value = "renewable energies' waste weighted at -1 so that there is no net loss in using renewable energies"

# Row 18
$ws.Range("A18").Value = "R6'"
$ws.Range("B18").Value = -2
$ws.Range("C18").Value = "nonrenewable energy waste"
$ws.Range("E18").Value = "nonrenewable energy waste is weighted higher than the weight of nonrenewable energy, to discourage fossil fuel use"

# Row 19
$ws.Range("A19").Value = "R18'"
$ws.Range("B19").Value = -1
$ws.Range("C19").Value = "metallic alloys waste"

# Row 20
$ws.Range("A20").Value = "R19'"
$ws.Range("B20").Value = -2
$ws.Range("C20").Value = "housing waste"

# Row 21
$ws.Range("A21").Value = "R20'"
$ws.Range("B21").Value = -1
$ws.Range("C21").Value = "electronics waste"

# Row 22
$ws.Range("A22").Value = "R21'"
$ws.Range("B22").Value = -1
$ws.Range("C22").Value = "farm waste"

# Row 23
$ws.Range("A23").Value = "R22'"
$ws.Range("B23").Value = -1
$ws.Range("C23").Value = "factory waste"

# --- Column widths (match the widened layout from the reorganized table)
$ws.Columns.Item(1).ColumnWidth = 8.498697916666666
$ws.Columns.Item(2).ColumnWidth = 5.666666666666667
$ws.Columns.Item(3).ColumnWidth = 22.498697916666668
$ws.Columns.Item(4).ColumnWidth = 14.330729166666666
$ws.Columns.Item(5).ColumnWidth = 159.49869791666666

# --- Restore the active selection to E10, matching the saved view state
$ws.Range("E10").Select() | Out-Null
